$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild data rows 2-13 (A:T) to reflect the updated NATMI output (per Dr Hou advice)
$data = @(
    ,@("ECs", "Ctf1", "Il6st", "ECs", 1.0, 0.3333333333333333, 0.1923443333333333, 0.577033, 0.06576397757907125, 0.06576397757907126, 3.0, 1.0, 44.017783, 132.053349, 0.253483683026081, 0.253483683026081, 8.466571125946334, 76.199140133517, 0.0166700952471876, 0.0166700952471876)
    ,@("ECs", "Ctf1", "Il6st", "FAPs", 1.0, 0.3333333333333333, 0.1923443333333333, 0.577033, 0.06576397757907125, 0.06576397757907126, 3.0, 1.0, 99.76728333333334, 299.30185, 0.5745264004968147, 0.5745264004968147, 19.18967160122778, 172.70704441105, 0.03778314132085703, 0.03778314132085704)
    ,@("ECs", "Ctf1", "Il6st", "M2", 1.0, 0.3333333333333333, 0.1923443333333333, 0.577033, 0.06576397757907125, 0.06576397757907126, 3.0, 1.0, 10.852095, 32.556285, 0.06249358376701795, 0.06249358376701795, 2.087338978045, 18.786050802405, 0.004109826641689979, 0.004109826641689981)
    ,@("ECs", "Ctf1", "Il6st", "sCs", 1.0, 0.3333333333333333, 0.1923443333333333, 0.577033, 0.06576397757907125, 0.06576397757907126, 3.0, 1.0, 19.01418566666667, 57.042557, 0.1094963327100864, 0.1094963327100864, 3.657270865931222, 32.915437793381, 0.007200914369336649, 0.007200914369336651)
    ,@("FAPs", "Ctf1", "Il6st", "ECs", 3.0, 1.0, 1.496366666666667, 4.4891, 0.5116190438852002, 0.5116190438852003, 3.0, 1.0, 44.017783, 132.053349, 0.253483683026081, 0.253483683026081, 65.86674322176667, 592.8006889958999, 0.1296870795503027, 0.1296870795503027)
    ,@("FAPs", "Ctf1", "Il6st", "FAPs", 3.0, 1.0, 1.496366666666667, 4.4891, 0.5116190438852002, 0.5116190438852003, 3.0, 1.0, 99.76728333333334, 299.30185, 0.5745264004968147, 0.5745264004968147, 149.2884372038889, 1343.595934835, 0.2939386477089859, 0.293938647708986)
    ,@("FAPs", "Ctf1", "Il6st", "M2", 3.0, 1.0, 1.496366666666667, 4.4891, 0.5116190438852002, 0.5116190438852003, 3.0, 1.0, 10.852095, 32.556285, 0.06249358376701795, 0.06249358376701795, 16.2387132215, 146.1484189935, 0.03197290757584139, 0.0319729075758414)
    ,@("FAPs", "Ctf1", "Il6st", "sCs", 3.0, 1.0, 1.496366666666667, 4.4891, 0.5116190438852002, 0.5116190438852003, 3.0, 1.0, 19.01418566666667, 57.042557, 0.1094963327100864, 0.1094963327100864, 28.45219362541111, 256.0697426287, 0.05602040905007018, 0.0560204090500702)
    ,@("sCs", "Ctf1", "Il6st", "ECs", 3.0, 1.0, 1.236056333333333, 3.708169, 0.4226169785357285, 0.4226169785357286, 3.0, 1.0, 44.017783, 132.053349, 0.253483683026081, 0.253483683026081, 54.40845945644233, 489.6761351079809, 0.1071265082285907, 0.1071265082285907)
    ,@("sCs", "Ctf1", "Il6st", "FAPs", 3.0, 1.0, 1.236056333333333, 3.708169, 0.4226169785357285, 0.4226169785357286, 3.0, 1.0, 99.76728333333334, 299.30185, 0.5745264004968147, 0.5745264004968147, 123.3179824236278, 1109.86184181265, 0.2428046114669717, 0.2428046114669717)
    ,@("sCs", "Ctf1", "Il6st", "M2", 3.0, 1.0, 1.236056333333333, 3.708169, 0.4226169785357285, 0.4226169785357286, 3.0, 1.0, 10.852095, 32.556285, 0.06249358376701795, 0.06249358376701795, 13.413800754685, 120.724206792165, 0.02641084954948657, 0.02641084954948658)
    ,@("sCs", "Ctf1", "Il6st", "sCs", 3.0, 1.0, 1.236056333333333, 3.708169, 0.4226169785357285, 0.4226169785357286, 3.0, 1.0, 19.01418566666667, 57.042557, 0.1094963327100864, 0.1094963327100864, 23.50260461645922, 211.523441548133, 0.04627500929067958, 0.04627500929067959)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}
